$d = $word.ActiveDocument
$result = $d.Content.Find.Execute("area_atuacao", $true, $false, $false, $false, $false, $true, 1, $false, "atividade", 2)
Write-Output "Find result: $result"
